$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header/values for the new custom email feature column (B)
$ws.Range("B1").Value = "test1"
$ws.Range("B2").Value = "test2"

# Update the active selection to match the new edit location
$ws.Range("B3").Select()
